$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (pushes existing row 16.. down to 17..,
# matching the author's edit: a new price record was added for this
# market/product sheet, with everything from the old row 16 onward
# shifting down by one row).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44971
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108002
$ws.Range("J16").Value = "Mango"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 7000
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 7000
$ws.Range("Q16").Value = "$/bandeja 4 kilos"
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 1750
$ws.Range("T16").Value = 4
